$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format temporarily so numeric-looking strings
# (e.g. "145.02") are stored as text instead of being auto-converted to numbers,
# matching the original inlineStr cell type. Revert the style afterwards so no
# cell ends up with a residual explicit style.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "62.929.73"
$ws.Range("E2").Value = "  +5.34%  "

$ws.Range("D3").Value = "3.113.02"
$ws.Range("E3").Value = "  +3.49%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "584.17"
$ws.Range("E5").Value = "  +3.45%  "

$ws.Range("D6").Value = "145.02"
$ws.Range("E6").Value = "  +3.30%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").Value = "3.109.69"
$ws.Range("E8").Value = "  +3.80%  "

$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  +1.76%  "

$ws.Range("D10").Value = "0.150"
$ws.Range("E10").Value = "  +12.00%  "

$ws.Range("E11").Value = "  +7.72%  "

$ws.Range("E12").Value = "  +1.89%  "

$ws.Range("D13").Value = "0.0000249"
$ws.Range("E13").Value = "  +7.68%  "

$ws.Range("D14").Value = "35.67"
$ws.Range("E14").Value = "  +5.13%  "

$ws.Range("E15").Value = "  +0.08%  "

$ws.Range("D16").Value = "3.627.20"
$ws.Range("E16").Value = "  +3.48%  "

$ws.Range("D17").Value = "7.18"
$ws.Range("E17").Value = "  -0.43%  "

$ws.Range("D18").Value = "3.114.67"
$ws.Range("E18").Value = "  +3.47%  "

$ws.Range("D19").Value = "62.820.15"
$ws.Range("E19").Value = "  +5.20%  "

$ws.Range("D20").Value = "469.16"
$ws.Range("E20").Value = "  +7.36%  "

$ws.Range("D21").Value = "14.11"
$ws.Range("E21").Value = "  +3.30%  "

$ws.Range("D22").Value = "0.729"
$ws.Range("E22").Value = "  +1.49%  "

$ws.Range("D23").Value = "7.55"
$ws.Range("E23").Value = "  +6.23%  "

$ws.Range("D24").Value = "13.32"
$ws.Range("E24").Value = "  -0.35%  "

$ws.Range("D25").Value = "82.03"
$ws.Range("E25").Value = "  +1.80%  "

$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("E27").Value = "  +1.16%  "

$ws.Range("E28").Value = "  +4.92%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "8.29"
$ws.Range("E29").Value = "  +6.08%  "

$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("D31").Value = "6.84"
$ws.Range("E31").Value = "  +7.98%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.111"
$ws.Range("E32").Value = "  +5.59%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "27.05"
$ws.Range("E33").Value = "  +4.23%  "

$ws.Range("D34").Value = "0.0₃0874"
$ws.Range("E34").Value = "  +10.80%  "

$ws.Range("D35").Value = "2.39"
$ws.Range("E35").Value = "  +13.53%  "

$ws.Range("E36").Value = "  +4.30%  "

$ws.Range("D37").Value = "6.06"
$ws.Range("E37").Value = "  +2.37%  "

$ws.Range("D38").Value = "3.29"
$ws.Range("E38").Value = "  +17.35%  "

$ws.Range("D39").Value = "51.03"
$ws.Range("E39").Value = "  +3.73%  "

$ws.Range("D40").Value = "432.45"
$ws.Range("E40").Value = "  +7.51%  "

$ws.Range("E41").Value = "  +2.07%  "

$ws.Range("D42").Value = "2.932.97"
$ws.Range("E42").Value = "  +6.11%  "

$ws.Range("E43").Value = "  +4.61%  "

$ws.Range("E44").Value = "  +10.38%  "

$ws.Range("D45").Value = "0.112"
$ws.Range("E45").Value = "  +4.30%  "

$ws.Range("E46").Value = "  +5.95%  "

$ws.Range("D47").Value = "35.49"
$ws.Range("E47").Value = "  +4.13%  "

$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("D49").Value = "123.66"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("E50").Value = "  +0.97%  "

$ws.Range("D51").Value = "24.65"
$ws.Range("E51").Value = "  +4.33%  "

# Restore the default (Normal) style on column D so no stray explicit
# cell style references remain after the text-format trick above.
$dRange.Style = "Normal"
